$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# --- Sheet1 ("CV Result") numeric updates rows 2-7, columns C-F ---
$ws1.Range("C2").Value = 0.9098967684021544
$ws1.Range("D2").Value = 0.8230952380952381
$ws1.Range("E2").Value = 0.6674521244763615
$ws1.Range("F2").Value = 0.5428571428571428
$ws1.Range("C3").Value = 0.8844703770197486
$ws1.Range("D3").Value = 0.7359523809523809
$ws1.Range("E3").Value = 0.6078994614003591
$ws1.Range("F3").Value = 0.4095238095238095
$ws1.Range("C4").Value = 0.9179136798088411
$ws1.Range("D4").Value = 0.7521582733812949
$ws1.Range("E4").Value = 0.6682795698924732
$ws1.Range("F4").Value = 0.4383693045563549
$ws1.Range("C5").Value = 0.8769489247311828
$ws1.Range("D5").Value = 0.8164268585131894
$ws1.Range("E5").Value = 0.5932795698924731
$ws1.Range("F5").Value = 0.5038369304556355
$ws1.Range("C6").Value = 0.8710722819593787
$ws1.Range("D6").Value = 0.8317745803357314
$ws1.Range("E6").Value = 0.5732078853046595
$ws1.Range("F6").Value = 0.5553956834532374
$ws1.Range("C7").Value = 0.892060406384261
$ws1.Range("D7").Value = 0.7887374461979915
$ws1.Range("E7").Value = 0.6220237221932653
$ws1.Range("F7").Value = 0.489996574169236

# --- Sheet2 ("95% Importance") updates rows 2-13 ---
$ws2.Range("B2").Value = 15
$ws2.Range("C2").Value = "status.of.existing.checking.account"
$ws2.Range("D2").Value = 471.2633125782013
$ws2.Range("E2").Value = 0.2709371718059145
$ws2.Range("F2").Value = 0.2709371718059145
$ws2.Range("B3").Value = 0
$ws2.Range("C3").Value = "duration.in.month"
$ws2.Range("D3").Value = 310.8317297816276
$ws2.Range("E3").Value = 0.178702367714227
$ws2.Range("F3").Value = 0.4496395395201415
$ws2.Range("B4").Value = 3
$ws2.Range("C4").Value = "age.in.years"
$ws2.Range("D4").Value = 167.268590426445
$ws2.Range("E4").Value = 0.0961655142942675
$ws2.Range("F4").Value = 0.545805053814409
$ws2.Range("B5").Value = 17
$ws2.Range("C5").Value = "purpose"
$ws2.Range("D5").Value = 127.706204354763
$ws2.Range("E5").Value = 0.07342043589316341
$ws2.Range("F5").Value = 0.6192254897075724
$ws2.Range("B6").Value = 6
$ws2.Range("C6").Value = "other.installment.plans"
$ws2.Range("D6").Value = 101.4248958945274
$ws2.Range("E6").Value = 0.05831087146172149
$ws2.Range("F6").Value = 0.6775363611692939
$ws2.Range("B7").Value = 10
$ws2.Range("C7").Value = "savings.account.and.bonds"
$ws2.Range("D7").Value = 101.0704591661692
$ws2.Range("E7").Value = 0.05810709984996549
$ws2.Range("F7").Value = 0.7356434610192595
$ws2.Range("B8").Value = 18
$ws2.Range("C8").Value = "credit.history"
$ws2.Range("D8").Value = 95.67135408818721
$ws2.Range("E8").Value = 0.05500306390855399
$ws2.Range("F8").Value = 0.7906465249278135
$ws2.Range("B9").Value = 1
$ws2.Range("C9").Value = "installment.rate.in.percentage.of.disposable.income"
$ws2.Range("D9").Value = 73.13659776449202
$ws2.Range("E9").Value = 0.04204745505312398
$ws2.Range("F9").Value = 0.8326939799809375
$ws2.Range("B10").Value = 9
$ws2.Range("C10").Value = "other.debtors.or.guarantors"
$ws2.Range("D10").Value = 64.2222418129444
$ws2.Range("E10").Value = 0.03692244250595535
$ws2.Range("F10").Value = 0.8696164224868929
$ws2.Range("B11").Value = 2
$ws2.Range("C11").Value = "present.residence.since"
$ws2.Range("D11").Value = 44.09301280975342
$ws2.Range("E11").Value = 0.02534981159836023
$ws2.Range("F11").Value = 0.8949662340852531
$ws2.Range("B12").Value = 13
$ws2.Range("C12").Value = "present.employment.since"
$ws2.Range("D12").Value = 41.33051123321057
$ws2.Range("E12").Value = 0.02376160317160369
$ws2.Range("F12").Value = 0.9187278372568568
$ws2.Range("B13").Value = 16
$ws2.Range("C13").Value = "housing"
$ws2.Range("D13").Value = 38.6174860715866
$ws2.Range("E13").Value = 0.02220183956448712
$ws2.Range("F13").Value = 0.9409296768213439

# --- Sheet3 ("Increase CV Filter") updates rows 2-8, then delete row 9 ---
$ws3.Range("B2").Value = 15
$ws3.Range("C2").Value = "status.of.existing.checking.account"
$ws3.Range("D2").Value = 471.2633125782013
$ws3.Range("E2").Value = 0.2709371718059145
$ws3.Range("F2").Value = 0.2709371718059145
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = "duration.in.month"
$ws3.Range("D3").Value = 310.8317297816276
$ws3.Range("E3").Value = 0.178702367714227
$ws3.Range("F3").Value = 0.4496395395201415
$ws3.Range("B4").Value = 17
$ws3.Range("C4").Value = "purpose"
$ws3.Range("D4").Value = 127.706204354763
$ws3.Range("E4").Value = 0.07342043589316341
$ws3.Range("F4").Value = 0.6192254897075724
$ws3.Range("B5").Value = 6
$ws3.Range("C5").Value = "other.installment.plans"
$ws3.Range("D5").Value = 101.4248958945274
$ws3.Range("E5").Value = 0.05831087146172149
$ws3.Range("F5").Value = 0.6775363611692939
$ws3.Range("B6").Value = 10
$ws3.Range("C6").Value = "savings.account.and.bonds"
$ws3.Range("D6").Value = 101.0704591661692
$ws3.Range("E6").Value = 0.05810709984996549
$ws3.Range("F6").Value = 0.7356434610192595
$ws3.Range("B7").Value = 18
$ws3.Range("C7").Value = "credit.history"
$ws3.Range("D7").Value = 95.67135408818721
$ws3.Range("E7").Value = 0.05500306390855399
$ws3.Range("F7").Value = 0.7906465249278135
$ws3.Range("B8").Value = 9
$ws3.Range("C8").Value = "other.debtors.or.guarantors"
$ws3.Range("D8").Value = 64.2222418129444
$ws3.Range("E8").Value = 0.03692244250595535
$ws3.Range("F8").Value = 0.8696164224868929
$ws3.Rows(9).Delete()

# --- Sheet4 ("Increase CV STEP") updates rows 2-13: C,D numeric + E text ---
$ws4.Range("C2").Value = 0.7079597464885234
$ws4.Range("D2").Value = 0.2079597464885234
$ws4.Range("E2").Value = "['status.of.existing.checking.account']"
$ws4.Range("C3").Value = 0.7553816375471051
$ws4.Range("D3").Value = 0.04742189105858163
$ws4.Range("E3").Value = "['duration.in.month']"
$ws4.Range("C4").Value = 0.7485474477560808
$ws4.Range("D4").Value = -0.006834189791024237
$ws4.Range("E4").Value = "['age.in.years']"
$ws4.Range("C5").Value = 0.7674050188420692
$ws4.Range("D5").Value = 0.01202338129496416
$ws4.Range("E5").Value = "['purpose']"
$ws4.Range("C6").Value = 0.7715508735868448
$ws4.Range("D6").Value = 0.004145854744775579
$ws4.Range("E6").Value = "['other.installment.plans']"
$ws4.Range("C7").Value = 0.7766170777663584
$ws4.Range("D7").Value = 0.005066204179513534
$ws4.Range("E7").Value = "['savings.account.and.bonds']"
$ws4.Range("C8").Value = 0.7838550017129153
$ws4.Range("D8").Value = 0.007237923946556979
$ws4.Range("E8").Value = "['credit.history']"
$ws4.Range("C9").Value = 0.7826848235697156
$ws4.Range("D9").Value = -0.001170178143199752
$ws4.Range("E9").Value = "['installment.rate.in.percentage.of.disposable.income']"
$ws4.Range("C10").Value = 0.7951287255909558
$ws4.Range("D10").Value = 0.01127372387804049
$ws4.Range("E10").Value = "['other.debtors.or.guarantors']"
$ws4.Range("C11").Value = 0.793486125385406
$ws4.Range("D11").Value = -0.001642600205549849
$ws4.Range("E11").Value = "['present.residence.since']"
$ws4.Range("C12").Value = 0.7958266529633435
$ws4.Range("D12").Value = 0.0006979273723877188
$ws4.Range("E12").Value = "['present.employment.since']"
$ws4.Range("C13").Value = 0.7943931997259335
$ws4.Range("D13").Value = -0.0007355258650223462
$ws4.Range("E13").Value = "['housing']"
